# Penalty Reward System (unfinished) - remove stale PO rows and
# refresh the dependent volume-insights summary.

$wb = $excel.ActiveWorkbook

# --- "Daily PO" sheet: drop the last two PO records (rows 4 & 5) ---
$wsDailyPo = $wb.Worksheets.Item("Daily PO")
$wsDailyPo.Rows.Item(4).EntireRow.Delete()
$wsDailyPo.Rows.Item(4).EntireRow.Delete()

# --- "Merged (Optional)" sheet: drop the corresponding merged rows
#     (the 2023-06-14 and 2023-07-12 entries, originally rows 7 & 8) ---
$wsMerged = $wb.Worksheets.Item("Merged (Optional)")
$wsMerged.Rows.Item(7).EntireRow.Delete()
$wsMerged.Rows.Item(7).EntireRow.Delete()

# --- "PO Volume Insights" sheet: recompute total/average PO quantity
#     now that only 2 PO rows remain (112 + 8 = 120, average 60) ---
$wsInsights = $wb.Worksheets.Item("PO Volume Insights")
$wsInsights.Range("A2").Value = 120
$wsInsights.Range("B2").Value = 60
